$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.08840041403804122
$ws.Cells.Item(1, 2).Value = 0.08838409956661764
$ws.Cells.Item(2, 1).Value = -0.07956733182590803
$ws.Cells.Item(2, 2).Value = 0.07945752881944301
$ws.Cells.Item(3, 1).Value = -0.029759923704434144
$ws.Cells.Item(3, 2).Value = 0.029552791897557995
$ws.Cells.Item(4, 1).Value = -0.02155279199026161
$ws.Cells.Item(4, 2).Value = 0.021266740098390358
$ws.Cells.Item(5, 1).Value = -0.01826674014111518
$ws.Cells.Item(5, 2).Value = 0.017303403014439667
$ws.Cells.Item(6, 1).Value = -0.003197793728906717
$ws.Cells.Item(6, 2).Value = 0.0031024489263415944
$ws.Cells.Item(7, 1).Value = 0.006897550954407805
$ws.Cells.Item(7, 2).Value = -0.006910435981331009
$ws.Cells.Item(8, 1).Value = 0.01691043586282337
$ws.Cells.Item(8, 2).Value = -0.016927540112898765
$ws.Cells.Item(9, 1).Value = 0.01892754007970643
$ws.Cells.Item(9, 2).Value = -0.018944560957309786
$ws.Cells.Item(10, 1).Value = 0.02094456092733843
$ws.Cells.Item(10, 2).Value = -0.02094434802821432
$ws.Cells.Item(11, 1).Value = 0.023944347988084758
$ws.Cells.Item(11, 2).Value = -0.023948043915844153
$ws.Cells.Item(12, 1).Value = 0.0042518052876063805
$ws.Cells.Item(12, 2).Value = -0.004662708312158426
$ws.Cells.Item(13, 1).Value = 0.008162708269657593
$ws.Cells.Item(13, 2).Value = -0.008378278410398288
$ws.Cells.Item(14, 1).Value = 0.01637827832243044
$ws.Cells.Item(14, 2).Value = -0.01649298352684525
$ws.Cells.Item(15, 1).Value = 0.017492983512977567
$ws.Cells.Item(15, 2).Value = -0.017590908327328236
$ws.Cells.Item(16, 1).Value = -0.006032695468067839
$ws.Cells.Item(16, 2).Value = 0.006002986800311305
$ws.Cells.Item(17, 1).Value = -0.004002986823491206
$ws.Cells.Item(17, 2).Value = 0.003999999955776268
$ws.Cells.Item(18, 1).Value = -0.0028167678033987897
$ws.Cells.Item(18, 2).Value = 0.0028065338939669005
$ws.Cells.Item(19, 1).Value = 0.0011934660638703853
$ws.Cells.Item(19, 2).Value = -0.0012660199773524283
$ws.Cells.Item(20, 1).Value = -0.008015651452087269
$ws.Cells.Item(20, 2).Value = 0.008005481947170878
$ws.Cells.Item(21, 1).Value = -0.004005481990044579
$ws.Cells.Item(21, 2).Value = 0.003999999956706191
$ws.Cells.Item(22, 1).Value = -0.04569760516423749
$ws.Cells.Item(22, 2).Value = 0.0454890507015655
$ws.Cells.Item(23, 1).Value = -0.04048905076448861
$ws.Cells.Item(23, 2).Value = 0.040096849119257705
$ws.Cells.Item(24, 1).Value = -0.0200968493444158
$ws.Cells.Item(24, 2).Value = 0.01999999977139133
$ws.Cells.Item(25, 1).Value = -0.012105609321917399
$ws.Cells.Item(25, 2).Value = 0.012055164227996684
$ws.Cells.Item(26, 1).Value = -0.009555164270015126
$ws.Cells.Item(26, 2).Value = 0.009491795083235388
$ws.Cells.Item(27, 1).Value = -0.006991795125785849
$ws.Cells.Item(27, 2).Value = 0.006626946360399177
$ws.Cells.Item(28, 1).Value = -0.004626946400157372
$ws.Cells.Item(28, 2).Value = 0.004395145909436238
$ws.Cells.Item(29, 1).Value = 0.0026048539980525476
$ws.Cells.Item(29, 2).Value = -0.0026626771880904343
$ws.Cells.Item(30, 1).Value = 0.06266267654892443
$ws.Cells.Item(30, 2).Value = -0.06300639443405709
$ws.Cells.Item(31, 1).Value = 0.008252328443143497
$ws.Cells.Item(31, 2).Value = -0.008281564381789508
$ws.Cells.Item(32, 1).Value = -0.004000885620975225
$ws.Cells.Item(32, 2).Value = 0.003999999951950883
